$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1446.4166
$ws.Range("J17").Value = 1532.5454
$ws.Range("L17").Value = 4597.6362
$ws.Range("N17").Value = -4933.6362

$ws.Range("H32").Value = 16670267
$ws.Range("J32").Value = 4999.6665
$ws.Range("L32").Value = 4999.6665
$ws.Range("N32").Value = -5651.6665

$ws.Range("H98").Value = 4435.6
$ws.Range("J98").Value = 5718.25
$ws.Range("L98").Value = 5718.25
$ws.Range("N98").Value = -8714.25

$ws.Range("H118").Value = 431
$ws.Range("I118").Value = 349.85715
$ws.Range("K118").Value = 1049.57145
$ws.Range("M118").Value = 607.4285500000001

$ws.Range("H122").Value = 4435.6
$ws.Range("J122").Value = 5718.25
$ws.Range("L122").Value = 17154.75
$ws.Range("N122").Value = -22054.75

$ws.Range("H137").Value = 6637.8667
$ws.Range("I137").Value = 2754.75
$ws.Range("J137").Value = 22170.334
$ws.Range("K137").Value = 8264.25
$ws.Range("L137").Value = 66511.00199999999
$ws.Range("M137").Value = -5714.25
$ws.Range("N137").Value = -71611.00199999999

$ws.Range("H141").Value = 1901.8846
$ws.Range("I141").Value = 1854.6364
$ws.Range("K141").Value = 5563.9092
$ws.Range("M141").Value = -383.9092000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 160836.38
$ws.Range("I74").Value = 243237.61
$ws.Range("K74").Value = 243237.61
$ws.Range("M74").Value = -242363.61

$ws.Range("H77").Value = 160836.38
$ws.Range("I77").Value = 243237.61
$ws.Range("K77").Value = 1216188.05
$ws.Range("M77").Value = -1211820.05

$ws.Range("H110").Value = 1378.579
$ws.Range("I110").Value = 1257
$ws.Range("K110").Value = 1257
$ws.Range("M110").Value = 788

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 104760.5
$ws.Range("J132").Value = 104760.5
$ws.Range("L132").Value = 104760.5
$ws.Range("N132").Value = -114880.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2250.625
$ws.Range("I16").Value = 2143.5715
$ws.Range("K16").Value = 2143.5715
$ws.Range("M16").Value = -1856.5715

$ws.Range("H31").Value = 2980618.8
$ws.Range("J31").Value = 5213376.5
$ws.Range("L31").Value = 5213376.5
$ws.Range("N31").Value = -5213966.5

$ws.Range("H34").Value = 2980618.8
$ws.Range("J34").Value = 5213376.5
$ws.Range("L34").Value = 5213376.5
$ws.Range("N34").Value = -5213780.5

$ws.Range("H58").Value = 2644.6365
$ws.Range("I58").Value = 1974.25
$ws.Range("J58").Value = 3027.7144
$ws.Range("K58").Value = 1974.25
$ws.Range("L58").Value = 3027.7144
$ws.Range("M58").Value = -1771.25
$ws.Range("N58").Value = -3433.7144

$ws.Range("H113").Value = 2250.625
$ws.Range("I113").Value = 2143.5715
$ws.Range("K113").Value = 2143.5715
$ws.Range("M113").Value = 26.42849999999999

$ws.Range("H136").Value = 2644.6365
$ws.Range("I136").Value = 1974.25
$ws.Range("J136").Value = 3027.7144
$ws.Range("K136").Value = 5922.75
$ws.Range("L136").Value = 9083.143199999999
$ws.Range("M136").Value = -3372.75
$ws.Range("N136").Value = -14183.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 473.3158
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H39").Value = 10767.2
$ws.Range("J39").Value = 12834
$ws.Range("L39").Value = 38502
$ws.Range("N39").Value = -39090

$ws.Range("H55").Value = 4850.0435
$ws.Range("J55").Value = 5714.7646
$ws.Range("L55").Value = 17144.2938
$ws.Range("N55").Value = -17498.2938

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5315.077
$ws.Range("I70").Value = 2998.3333
$ws.Range("J70").Value = 6010.1
$ws.Range("K70").Value = 2998.3333
$ws.Range("L70").Value = 6010.1
$ws.Range("M70").Value = -2728.3333
$ws.Range("N70").Value = -6550.1

$ws.Range("H73").Value = 5315.077
$ws.Range("I73").Value = 2998.3333
$ws.Range("J73").Value = 6010.1
$ws.Range("K73").Value = 2998.3333
$ws.Range("L73").Value = 6010.1
$ws.Range("M73").Value = -2062.3333
$ws.Range("N73").Value = -7882.1

$ws.Range("H80").Value = 90911760
$ws.Range("I80").Value = 166668880
$ws.Range("K80").Value = 166668880
$ws.Range("M80").Value = -166667882

$ws.Range("H83").Value = 90911760
$ws.Range("I83").Value = 166668880
$ws.Range("K83").Value = 833344400
$ws.Range("M83").Value = -833339408

$ws.Range("H107").Value = 11425.777
$ws.Range("J107").Value = 17698.4
$ws.Range("L107").Value = 17698.4
$ws.Range("N107").Value = -21538.4

$ws.Range("H113").Value = 2263.3333
$ws.Range("I113").Value = 2034.3334
$ws.Range("K113").Value = 2034.3334
$ws.Range("M113").Value = 135.6666

$ws.Range("H126").Value = 7900.467
$ws.Range("I126").Value = 2130
$ws.Range("K126").Value = 6390
$ws.Range("M126").Value = -3920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1830
$ws.Range("N16").ClearContents()

$ws.Range("H22").Value = 1459.4615
$ws.Range("J22").Value = 749.5
$ws.Range("L22").Value = 749.5
$ws.Range("N22").Value = -1339.5

$ws.Range("H27").Value = 1459.4615
$ws.Range("J27").Value = 749.5
$ws.Range("L27").Value = 749.5
$ws.Range("N27").Value = -963.5

$ws.Range("H40").Value = 12703.385
$ws.Range("I40").Value = 24832.834
$ws.Range("K40").Value = 24832.834
$ws.Range("M40").Value = -24696.834

$ws.Range("H61").Value = 1549.5
$ws.Range("I61").Value = 1549.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1549.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1347.5
$ws.Range("N61").ClearContents()

$ws.Range("H68").Value = 3398.125
$ws.Range("I68").Value = 3383.5715
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 3383.5715
$ws.Range("L68").Value = 3500
$ws.Range("M68").Value = -2634.5715
$ws.Range("N68").Value = -4998

$ws.Range("H71").Value = 3398.125
$ws.Range("I71").Value = 3383.5715
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 16917.8575
$ws.Range("L71").Value = 17500
$ws.Range("M71").Value = -13173.8575
$ws.Range("N71").Value = -24988

$ws.Range("H113").Value = 1549.5
$ws.Range("I113").Value = 1549.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1549.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 620.5
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 4999.1
$ws.Range("I132").Value = 3880.6667
$ws.Range("K132").Value = 11642.0001
$ws.Range("M132").Value = -9112.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H62").Value = 5892
$ws.Range("I62").Value = 4648.8
$ws.Range("K62").Value = 4648.8
$ws.Range("M62").Value = -4024.8

$ws.Range("H65").Value = 5892
$ws.Range("I65").Value = 4648.8
$ws.Range("K65").Value = 23244
$ws.Range("M65").Value = -20124

$ws.Range("H132").Value = 1319.7949
$ws.Range("I132").Value = 1239.0968
$ws.Range("K132").Value = 3717.2904
$ws.Range("M132").Value = -1187.2904

$ws.Range("H136").Value = 4052.3547
$ws.Range("I136").Value = 2090.276
$ws.Range("K136").Value = 6270.828
$ws.Range("M136").Value = -3720.828
